# Update annotations for Ruilin
# 1. Fix B7 so it is stored as a real number (it was previously stored as text "3").
# 2. Append a new annotation row (row 8) with the additional Ruilin comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix B7: convert the text "3" into a genuine numeric value ---
$ws.Range("B7").Value = 3

# --- Add new row 8 ---
$ws.Range("A8").Value = "Ruilin"

# B8 must stay a text value ("2"), not be auto-converted to a number like B2..B6.
# Temporarily force a text number format, assign the value, then clear the
# formatting again so the cell ends up with no explicit style (matching the
# rest of the sheet) while still keeping its value stored as text.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2"
$ws.Range("B8").ClearFormats()

$ws.Range("C8").Value = "the language is convoluted"
$ws.Range("D8").Value = "CRT"
$ws.Range("E8").Value = "WRI"
$ws.Range("F8").Value = "f5b44bd7-9311-4cfc-b939-3b86c20706ac"
$ws.Range("G8").Value = "SkYXvCR6W_annotated.xlsx"
$ws.Range("H8").Value = "On top of this, I do not enjoy the style the paper is written in, the language is convoluted."
